$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Latest HO Xliff Generate Date / Correspond Handoff Datetime for file 2a73fa5a on Overview & de-de (shared string)
$wsOverview.Range("G4").Value = "2016-08-14 03:03:59"
$wsDeDe.Range("H4").Value = "2016-08-14 03:03:59"

# zh-cn row for file 2a73fa5a: Correspond Handoff Datetime (H4) and Correspond Handback DateTime (K4)
$wsZhCn.Range("H4").Value = "2016-08-14 03:03:50"
$wsZhCn.Range("K4").Value = "2016-08-14 03:04:20"

# de-de row for file 2a73fa5a: Correspond Handback DateTime (K4)
$wsDeDe.Range("K4").Value = "2016-08-14 03:04:31"
